$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name (B1 on both sheets) to the new value
$newProductName = "2619-MS-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-FEE-%INT-MORE-1st"
$ws1.Range("B1").Value = $newProductName
$ws2.Range("B1").Value = $newProductName

# Update short name (B2 on input sheet) from numeric 2619 to text "261x"
$ws1.Range("B2").Value = "261x"

# Update the selection on the input sheet to B3
$ws1.Range("B3").Select() | Out-Null

# Make the output sheet the active / selected tab
$ws2.Select() | Out-Null
$ws2.Range("B1").Select() | Out-Null
